# Inserts a new weekly price record for "Alcachofa" (Madrigal / Primera) at
# Macroferia Regional de Talca. The new row is inserted above row 48, which
# pushes the existing rows 48-82 down to 49-83 (dimension grows from R82 to R83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 48, shifting rows 48:82 down to 49:83.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A48").Value = 5
$ws.Range("B48").Value = "Macroferia Regional de Talca"
$ws.Range("C48").Value = "Maule"
$ws.Range("D48").Value = 44767
$ws.Range("E48").Value = 7
$ws.Range("F48").Value = 100112013
$ws.Range("G48").Value = "Alcachofa"
$ws.Range("H48").Value = "Madrigal"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 400
$ws.Range("K48").Value = 12000
$ws.Range("L48").Value = 12000
$ws.Range("M48").Value = 12000
$ws.Range("N48").Value = "$/caja 40 unidades"
$ws.Range("O48").Value = "Provincia del Elquí"
$ws.Range("P48").Value = 300
$ws.Range("Q48").Value = 40
$ws.Range("R48").Value = "Hortaliza"
